# Client form update:
#  - enforce parental relationship options (ref1_parentesco / ref2_parentesco)
#  - limit Fuente to LUZWARE/LEADS/SEGUIMIENTO, add fuente_base for LUZWARE
#  - add Usuario Cipre & Contrasena fields (editor-only; not part of this data row)
#  - store fuente_base; expose fields in editor
#  - replace the old single integration-test client row with two new sample rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old sample/test data row (row 2) entirely first.
$ws.Range("A2:AC2").ClearContents()

# Columns whose sample values look numeric/date-like to Excel's auto-detection
# (digit-only phone numbers, amounts, a plazo count, and an ISO date) need to be
# pre-formatted as Text so they are stored as shared strings, matching how the
# source data file already represents them (e.g. id/telefono/monto/plazo/fecha
# are free-text fields, not numeric cells).
$textCols = @("G2","I2","K2","L2","Q2","T2","X2","G3","I3","K3","L3","Q3","T3","X3")
foreach ($addr in $textCols) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2: Prueba Uno
$ws.Range("A2").Value = "C9998"
$ws.Range("B2").Value = "Nuevo"
$ws.Range("C2").Value = "MEJORAVIT"
$ws.Range("D2").Value = "LEADS"
$ws.Range("F2").Value = "Prueba Uno"
$ws.Range("G2").Value = "555000111"
$ws.Range("H2").Value = "BeworkEdoMex"
$ws.Range("I2").Value = "1000"
$ws.Range("K2").Value = "5000"
$ws.Range("L2").Value = "12"
$ws.Range("M2").Value = "Soltero"
$ws.Range("N2").Value = "Propia"
$ws.Range("O2").Value = "uno@example.com"
$ws.Range("P2").Value = "Ref Uno"
$ws.Range("Q2").Value = "111111"
$ws.Range("R2").Value = "Hijo"
$ws.Range("S2").Value = "Ref Dos"
$ws.Range("T2").Value = "222222"
$ws.Range("U2").Value = "Amig@"
$ws.Range("V2").Value = "2 años"
$ws.Range("X2").Value = "2025-12-11"
$ws.Range("Z2").Value = "DISPERSADO"

# Row 3: Prueba Dos
$ws.Range("A3").Value = "C9999"
$ws.Range("B3").Value = "Nuevo"
$ws.Range("C3").Value = "MEJORAVIT"
$ws.Range("D3").Value = "LEADS"
$ws.Range("F3").Value = "Prueba Dos"
$ws.Range("G3").Value = "555000222"
$ws.Range("H3").Value = "BeworkEdoMex"
$ws.Range("I3").Value = "1000"
$ws.Range("K3").Value = "5000"
$ws.Range("L3").Value = "12"
$ws.Range("M3").Value = "Soltero"
$ws.Range("N3").Value = "Propia"
$ws.Range("O3").Value = "uno@example.com"
$ws.Range("P3").Value = "Ref Uno"
$ws.Range("Q3").Value = "111111"
$ws.Range("R3").Value = "Hijo"
$ws.Range("S3").Value = "Ref Dos"
$ws.Range("T3").Value = "222222"
$ws.Range("U3").Value = "Amig@"
$ws.Range("V3").Value = "2 años"
$ws.Range("X3").Value = "2025-12-11"
$ws.Range("Z3").Value = "DISPERSADO"

# Restore the default (unstyled) cell style now that the values are locked in
# as text, so the sheet doesn't pick up new styling compared to the source.
foreach ($addr in $textCols) {
    $ws.Range($addr).Style = "Normal"
}
